# KIBON-1710 change header names
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E6").Value = "{bgMitSelbstbehaltTitel}"
$ws.Range("J6").Value = "{bgOhneSelbstbehaltTitel}"
